$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.011.48"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").Value = "3.072.56"
$ws.Range("E3").Value = "  +1.77%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "3.067.45"
$ws.Range("E8").Value = "  +1.65%  "

$ws.Range("E9").Value = "  +2.56%  "

$ws.Range("E10").Value = "  +3.39%  "

$ws.Range("E11").Value = "  +2.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("E13").Value = "  +6.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = "3.568.01"
$ws.Range("E15").Value = "  +1.56%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "62.997.19"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.113"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.38%  "

$ws.Range("D18").Value = "3.071.74"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.697"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("E24").Value = "  +2.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.44%  "

$ws.Range("E26").Value = "  +0.86%  "

$ws.Range("E27").Value = "  +1.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.56%  "

$ws.Range("B31").Value = "Mantle"
$ws.Range("C31").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.50%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.40%  "

$ws.Range("E33").Value = "  -0.42%  "

$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.12%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.15%  "

$ws.Range("E36").Value = "  +2.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "486.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").Value = "3.251.74"
$ws.Range("E38").Value = "  +4.94%  "

$ws.Range("E39").Value = "  +3.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0794"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.18%  "

$ws.Range("E41").Value = "  +2.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.73%  "

$ws.Range("E44").Value = "  +2.16%  "


$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.86%  "

$ws.Range("E48").Value = "  +1.17%  "

$ws.Range("E49").Value = "  +3.22%  "

$ws.Range("D50").Value = "0.0₃0522"
$ws.Range("E50").Value = "  +6.53%  "

$ws.Range("E51").Value = "  +2.70%  "

